$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.014.73"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.880.68"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "243.31"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.4967"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("D8").Value = "0.2929"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.06659"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "1.880.64"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "16.83"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "0.07251"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "0.6684"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("D14").Value = "86.69"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "4.928"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "29.986.72"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "0.000007923"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").Value = "0.9981"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "12.82"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").Value = "2.124.09"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "0.9987"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "4.789"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").Value = "5.778"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "9.101"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "143.23"
$ws.Range("E25").Value = "  +5.96%  "
$ws.Range("D26").Value = "149.92"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").Value = "17.12"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "1.923"
$ws.Range("E28").Value = "  -3.48%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "4.208"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "0.08781"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "3.978"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "0.05089"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7168"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "2.666"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "0.01827"
$ws.Range("E37").Value = "  +7.91%  "
$ws.Range("D38").Value = "2.693"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("D39").Value = "2.185"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("D40").Value = "0.9337"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.4274"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.805"
$ws.Range("E42").Value = "  -5.42%  "
$ws.Range("D44").Value = "102.54"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").Value = "7.474"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").Value = "0.1272"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "0.05665"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "32.70"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "0.3811"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.294"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "56.07"
$ws.Range("E51").Value = "  -1.17%  "
